# cryptos.xlsx -- "Updated symbol list" data refresh
# (Tue Jan 31 07:43:42 UTC 2023 with GitHub Actions)
#
# Refreshes the scraped coinranking.com price/volume snapshot: most rows
# just get new Price (D) / Volume 1h% (E) readings, two rows (20 & 21)
# swap which coin (MCDex / ProBitToken) occupies them along with their
# Coin/Link/Price/Volume columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E are stored as plain text in the workbook (e.g. "0.005800"
# or "-0.73%") even though they look numeric -- a bare assignment would let
# the host auto-convert them to real numbers and silently lose information
# (trailing zeros, the literal '%' sign, thousands separators, ...). Prefixing
# the value with a leading apostrophe forces it to be entered as text, exactly
# like typing '0.005800 into Excel.

$ws.Range("D2").Value = "'311.03"
$ws.Range("E2").Value = "'-0.73%"
$ws.Range("D3").Value = "'37.67"
$ws.Range("E3").Value = "'-3.96%"
$ws.Range("D4").Value = "'5.087"
$ws.Range("E4").Value = "'-0.85%"
$ws.Range("D5").Value = "'0.07773"
$ws.Range("E5").Value = "'-4.35%"
$ws.Range("D6").Value = "'4.347"
$ws.Range("E6").Value = "'-3.21%"
$ws.Range("D7").Value = "'1.899"
$ws.Range("E7").Value = "'-3.18%"
$ws.Range("D8").Value = "'8.208"
$ws.Range("E8").Value = "'-1.00%"
$ws.Range("E9").Value = "'-7.45%"
$ws.Range("D10").Value = "'0.9162"
$ws.Range("E10").Value = "'-2.49%"
$ws.Range("E11").Value = "'-8.98%"
$ws.Range("D12").Value = "'0.1915"
$ws.Range("E12").Value = "'-2.56%"
$ws.Range("D13").Value = "'0.09271"
$ws.Range("E13").Value = "'3.11%"
$ws.Range("D14").Value = "'0.03406"
$ws.Range("E14").Value = "'-2.44%"
$ws.Range("D15").Value = "'0.09691"
$ws.Range("E15").Value = "'-0.21%"
$ws.Range("D16").Value = "'0.001363"
$ws.Range("E16").Value = "'-3.66%"
$ws.Range("D17").Value = "'0.005800"
$ws.Range("E17").Value = "'-7.49%"
$ws.Range("D18").Value = "'3.557"
$ws.Range("E18").Value = "'-0.60%"
$ws.Range("D19").Value = "'0.3374"
$ws.Range("E19").Value = "'-2.63%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'5.034"
$ws.Range("E20").Value = "'0.48%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1267"
$ws.Range("E21").Value = "'-2.55%"
$ws.Range("D22").Value = "'0.2588"
$ws.Range("E22").Value = "'3.87%"
$ws.Range("D23").Value = "'0.02102"
$ws.Range("E23").Value = "'5,583.20%"
$ws.Range("E24").Value = "'0.23%"
$ws.Range("D25").Value = "'0.001212"
$ws.Range("E25").Value = "'-2.76%"
$ws.Range("D26").Value = "'0.004256"
$ws.Range("E26").Value = "'-10.04%"
$ws.Range("E27").Value = "'-66.65%"
$ws.Range("D39").Value = "'0.02117"
$ws.Range("E39").Value = "'-4.39%"
$ws.Range("D40").Value = "'0.04950"
$ws.Range("E40").Value = "'-5.39%"
$ws.Range("D41").Value = "'0.007643"
$ws.Range("E41").Value = "'0.26%"
$ws.Range("D42").Value = "'0.009911"
$ws.Range("E42").Value = "'-4.37%"
$ws.Range("E43").Value = "'-3.79%"
$ws.Range("D44").Value = "'0.002059"
$ws.Range("E44").Value = "'-2.11%"
$ws.Range("D45").Value = "'0.008796"
$ws.Range("E45").Value = "'-3.63%"
$ws.Range("D46").Value = "'0.00006666"
$ws.Range("E46").Value = "'-2.15%"
$ws.Range("E47").Value = "'-0.21%"
$ws.Range("E48").Value = "'0.78%"
$ws.Range("E50").Value = "'-0.21%"
$ws.Range("E51").Value = "'-0.21%"
